$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("ProductionRates")

# 1. Insert the new worksheet right after "ProductionRates"
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $src)
$new.Name = "PadRates"

# 2. Copy cell formatting from "ProductionRates" (which has a near identical
#    layout, minus the Tank-id column) so the new sheet picks up the exact
#    same styles used throughout the workbook.

# Title cell
$src.Range("A1").Copy()
$new.Range("A1").PasteSpecial(-4122)

# Header row (blank corner cell + week labels T1..T5)
$src.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)
$src.Range("C2:G2").Copy()
$new.Range("B2:F2").PasteSpecial(-4122)

# Data rows - one row per pad (first occurrence of each pad id in the
# source sheet carries the right "interior" row styling).
$src.Range("A3:A3").Copy()
$new.Range("A3:A3").PasteSpecial(-4122)
$src.Range("C3:G3").Copy()
$new.Range("B3:F3").PasteSpecial(-4122)

$src.Range("A5:A5").Copy()
$new.Range("A4:A4").PasteSpecial(-4122)
$src.Range("C5:G5").Copy()
$new.Range("B4:F4").PasteSpecial(-4122)

$src.Range("A8:A8").Copy()
$new.Range("A5:A5").PasteSpecial(-4122)
$src.Range("C8:G8").Copy()
$new.Range("B5:F5").PasteSpecial(-4122)

$src.Range("A12:A12").Copy()
$new.Range("A6:A6").PasteSpecial(-4122)
$src.Range("C12:G12").Copy()
$new.Range("B6:F6").PasteSpecial(-4122)

# Last data row (bottom-bordered / thick-bottom styling)
$src.Range("A16:A16").Copy()
$new.Range("A7:A7").PasteSpecial(-4122)
$src.Range("C16:G16").Copy()
$new.Range("B7:F7").PasteSpecial(-4122)

# Trailing decorative formatting block (rows 19-27, cols E:G in the source
# become rows 10-18, cols D:F in the new sheet)
$src.Range("E19:G27").Copy()
$new.Range("D10:F18").PasteSpecial(-4122)

# 3. Fill in the actual values
$new.Range("A1").Value = "Table of Production Rate Forecasts by Pads [bbl/day]"

$new.Range("B2").Value = "T1"
$new.Range("C2").Value = "T2"
$new.Range("D2").Value = "T3"
$new.Range("E2").Value = "T4"
$new.Range("F2").Value = "T5"

$new.Range("A3").Value = "PP01"
$new.Range("B3").Value = 2116
$new.Range("C3").Value = 2058
$new.Range("D3").Value = 1998
$new.Range("E3").Value = 1996
$new.Range("F3").Value = 1992

$new.Range("A4").Value = "PP02"
$new.Range("B4").Value = 1398
$new.Range("C4").Value = 1380
$new.Range("D4").Value = 1374
$new.Range("E4").Value = 1371
$new.Range("F4").Value = 1365

$new.Range("A5").Value = "PP03"
$new.Range("B5").Value = 800
$new.Range("C5").Value = 796
$new.Range("D5").Value = 792
$new.Range("E5").Value = 784
$new.Range("F5").Value = 780

$new.Range("A6").Value = "PP04"
$new.Range("B6").Value = 993
$new.Range("C6").Value = 990
$new.Range("D6").Value = 990
$new.Range("E6").Value = 987
$new.Range("F6").Value = 987

$new.Range("A7").Value = "PP05"
$new.Range("B7").Value = 1790
$new.Range("C7").Value = 1776
$new.Range("D7").Value = 1774
$new.Range("E7").Value = 1770
$new.Range("F7").Value = 1766

# 4. Sheet view / selection details
$new.Range("A2").Select()
$new.Activate()

Write-Host "PadRates sheet created"
